$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Title: drop the "[information only]" suffix
$meta.Range("B5").Value = "NG-Imm-A RegisterFacility"

# Date: bump to the new generation timestamp
$meta.Range("B8").Value = "2025-08-03T02:35:31+01:00"

# --- Elements sheet --------------------------------------------------------
$els = $wb.Worksheets.Item("Elements")

# facilityStatus.Type(s) changed from "string" to "boolean"
$els.Range("K8").Value = "boolean`n"

# Append two new leaf elements: longitude (row 45) and latitude (row 46),
# mirroring the formatting of the last existing data row (44).
$els.Range("A44:AJ44").Copy()
$els.Range("A45:AJ46").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 45 - register-facility-model.longitude
$els.Range("A45").Value = "register-facility-model.longitude"
$els.Range("B45").Value = "register-facility-model.longitude"
$els.Range("D45").Value = ""
$els.Range("F45").Value = "0"
$els.Range("G45").Value = "1"
$els.Range("H45").Value = ""
$els.Range("I45").Value = ""
$els.Range("J45").Value = ""
$els.Range("K45").Value = "string`n"
$els.Range("L45").Value = "The Longitude of the health facility"
$els.Range("M45").Value = "The Longitude of the health facility"
$els.Range("P45").Value = ""
$els.Range("R45").Value = ""
$els.Range("S45").Value = ""
$els.Range("T45").Value = ""
$els.Range("U45").Value = ""
$els.Range("V45").Value = ""
$els.Range("W45").Value = ""
$els.Range("X45").Value = ""
$els.Range("Y45").Value = ""
$els.Range("Z45").Value = ""
$els.Range("AA45").Value = ""
$els.Range("AB45").Value = ""
$els.Range("AC45").Value = ""
$els.Range("AD45").Value = ""
$els.Range("AE45").Value = ""
$els.Range("AF45").Value = "register-facility-model.longitude"
$els.Range("AG45").Value = "0"
$els.Range("AH45").Value = "1"
$els.Range("AI45").Value = ""
$els.Range("AJ45").Value = ""

# Row 46 - register-facility-model.latitude
$els.Range("A46").Value = "register-facility-model.latitude"
$els.Range("B46").Value = "register-facility-model.latitude"
$els.Range("D46").Value = ""
$els.Range("F46").Value = "0"
$els.Range("G46").Value = "1"
$els.Range("H46").Value = ""
$els.Range("I46").Value = ""
$els.Range("J46").Value = ""
$els.Range("K46").Value = "string`n"
$els.Range("L46").Value = "The Latitude of the health facility"
$els.Range("M46").Value = "The Latitude of the health facility"
$els.Range("P46").Value = ""
$els.Range("R46").Value = ""
$els.Range("S46").Value = ""
$els.Range("T46").Value = ""
$els.Range("U46").Value = ""
$els.Range("V46").Value = ""
$els.Range("W46").Value = ""
$els.Range("X46").Value = ""
$els.Range("Y46").Value = ""
$els.Range("Z46").Value = ""
$els.Range("AA46").Value = ""
$els.Range("AB46").Value = ""
$els.Range("AC46").Value = ""
$els.Range("AD46").Value = ""
$els.Range("AE46").Value = ""
$els.Range("AF46").Value = "register-facility-model.latitude"
$els.Range("AG46").Value = "0"
$els.Range("AH46").Value = "1"
$els.Range("AI46").Value = ""
$els.Range("AJ46").Value = ""

Write-Host "done"
